$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from the last existing data row (59) down to
# the two new rows (60 and 61) so the new rows match the workbook's
# established per-column styling (bold/border/centered index column,
# date-formatted match-date column).
$ws.Range("A59:V59").Copy()
$ws.Range("A60:V60").PasteSpecial(-4122)
$ws.Range("A61:V61").PasteSpecial(-4122)

# Row 60: Alashkert 1 - 0 BKMA
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "armenia"
$ws.Range("C60").Value = "premier-league"
$ws.Range("D60").Value = "2023-2024"
$ws.Range("E60").Value = 45224.54166666666
$ws.Range("F60").Value = "Alashkert"
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = "BKMA"
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1.33
$ws.Range("K60").Value = "24/10/2023 00:12"
$ws.Range("L60").Value = 1.36
$ws.Range("M60").Value = "25/10/2023 12:51"
$ws.Range("N60").Value = 4.78
$ws.Range("O60").Value = "24/10/2023 00:12"
$ws.Range("P60").Value = 4.99
$ws.Range("Q60").Value = "25/10/2023 12:51"
$ws.Range("R60").Value = 7.01
$ws.Range("S60").Value = "24/10/2023 00:12"
$ws.Range("T60").Value = 8.109999999999999
$ws.Range("U60").Value = "25/10/2023 12:51"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/armenia/premier-league/alashkert-bkma/0Czsnty8/"

# Row 61: Pyunik Yerevan 3 - 1 Urartu
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "armenia"
$ws.Range("C61").Value = "premier-league"
$ws.Range("D61").Value = "2023-2024"
$ws.Range("E61").Value = 45224.54166666666
$ws.Range("F61").Value = "Pyunik Yerevan"
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = "Urartu"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1.66
$ws.Range("K61").Value = "24/10/2023 03:12"
$ws.Range("L61").Value = 1.61
$ws.Range("M61").Value = "25/10/2023 12:58"
$ws.Range("N61").Value = 3.7
$ws.Range("O61").Value = "24/10/2023 03:12"
$ws.Range("P61").Value = 3.92
$ws.Range("Q61").Value = "25/10/2023 12:58"
$ws.Range("R61").Value = 4.36
$ws.Range("S61").Value = "24/10/2023 03:12"
$ws.Range("T61").Value = 5.58
$ws.Range("U61").Value = "25/10/2023 12:58"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-urartu/tCTUlrik/"
